# Fruta / hortaliza, semanal
# Weekly refresh: reassigns the Fecha/Volumen/Precio values across the
# existing data rows (rows are re-shuffled onto a new weekly cycle while
# the market/product descriptive columns stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45084
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18500
$ws.Range("S2").Value = 1028

# Row 4
$ws.Range("D4").Value = 45072
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("S4").Value = 833

# Row 5
$ws.Range("D5").Value = 45092
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 19000
$ws.Range("P5").Value = 18667
$ws.Range("S5").Value = 1037

# Row 7
$ws.Range("D7").Value = 45083
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("S7").Value = 833

# Row 8
$ws.Range("D8").Value = 45055
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("S8").Value = 833

# Row 9
$ws.Range("D9").Value = 45076
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("S9").Value = 833

# Row 10
$ws.Range("D10").Value = 45085
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 19000
$ws.Range("O10").Value = 19000
$ws.Range("P10").Value = 19000
$ws.Range("S10").Value = 1056

# Row 11
$ws.Range("D11").Value = 45069
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("S11").Value = 833
